$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Daily auto-push: a new reading for 2026/02/01 (日, hour 13) was appended to
# the log. In the sheet it slots in right after the existing 2026/02/01 rows
# (749, 750) and before the 2026/12/29 block, so every row from the old 751
# through the old 792 shifts down by one (751 -> 752 ... 792 -> 793).
$ws.Rows(751).Insert()

# A751 ("2026/02/01") must stay literal text, not get auto-parsed into a
# date serial by Excel's input parser. Format as text while assigning, then
# clear the formatting again so the new row ends up styled like its
# neighbours (no explicit number format / style).
$ws.Range("A751").NumberFormat = "@"
$ws.Range("A751").Value = "2026/02/01"
$ws.Range("A751").ClearFormats()

$ws.Range("B751").Value = "日"
$ws.Range("C751").Value = 13
$ws.Range("D751").Value = 201
